$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header columns for the season record (Wins, Losses, Ties),
# matching the header style used by the other header cells (e.g. A1)
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in the season record for every data row (2-47) with the team's
# Wins/Losses/Ties totals
$lastRow = 47
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 63
    $ws.Cells.Item($r, 31).Value = 99
    $ws.Cells.Item($r, 32).Value = 0
}
